$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.594.68"
$ws.Range("D3").Value = "3.023.81"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.14"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.72"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.022.64"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("E14").Value = "  -4.95%  "
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "3.523.92"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "62.598.86"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "3.024.99"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.91"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -4.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.18"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.40"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.39"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "0.0₃0790"
$ws.Range("E36").Value = "  -3.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.58"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.07"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -9.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "420.28"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").Value = "2.782.07"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0353"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.44"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.39"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.08"
$ws.Range("E51").Value = "  -3.11%  "
